# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 4
    4  = 1
    5  = 3
    6  = 1
    7  = 2
    8  = 0
    9  = 2
    10 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
